$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.590.55"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "3.390.17"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'574.49"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").Value = "'137.69"
$ws.Range("E6").Value = "  +7.64%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.389.04"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "'0.477"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'7.47"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +7.96%  "
$ws.Range("D12").Value = "'0.394"
$ws.Range("E12").Value = "  +6.54%  "
$ws.Range("D13").Value = "3.971.96"
$ws.Range("E13").Value = "  +3.35%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("E15").Value = "  +7.30%  "
$ws.Range("D16").Value = "3.385.21"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "'25.32"
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").Value = "61.675.55"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("D19").Value = "'14.05"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").Value = "'5.89"
$ws.Range("E20").Value = "  +4.65%  "
$ws.Range("D21").Value = "'9.37"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "'388.28"
$ws.Range("E22").Value = "  +10.45%  "
$ws.Range("D23").Value = "'0.571"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("D24").Value = "3.530.47"
$ws.Range("E24").Value = "  +3.32%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000127"
$ws.Range("E25").Value = "  +17.88%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'71.03"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'1.60"
$ws.Range("E28").Value = "  +11.95%  "
$ws.Range("D29").Value = "'7.68"
$ws.Range("E29").Value = "  +6.73%  "
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "  +6.26%  "
$ws.Range("E32").Value = "  +5.87%  "
$ws.Range("D33").Value = "'2.15"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D35").Value = "3.423.92"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").Value = "'23.46"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("D37").Value = "'5.49"
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("D40").Value = "'162.37"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "'0.0795"
$ws.Range("E41").Value = "  +6.51%  "
$ws.Range("E42").Value = "  +12.76%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "'1.22"
$ws.Range("E44").Value = "  +7.63%  "
$ws.Range("D45").Value = "'0.772"
$ws.Range("E45").Value = "  +4.53%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'24.61"
$ws.Range("E48").Value = "  +8.74%  "
$ws.Range("D49").Value = "'6.97"
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("D50").Value = "'22.94"
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "2.372.81"
$ws.Range("E51").Value = "  +9.97%  "
